$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New table data (Collector / Team / Cycle / ... ) replacing the previous
# 14-row table with a fresh 13-row table. Columns E (Repayment_amount) and
# F (Pending Amount) hold text-looking numbers ("1,234.00") that must stay
# literal text, so we flip those cells to the Text number format before
# assigning the value (otherwise the engine parses the comma-grouped digits
# as a real number).
# ---------------------------------------------------------------------------

$rows = @(
    @("Shofa Khairunnisa",        "Anisa_s1", "s1", 1, "194,170.00",   "0.00"),
    @("Febri Fransiska",          "Anisa_s1", "s1", 6, "1,291,385.00", "0.00"),
    @("Raina Claresta Purwanika", "Anisa_s1", "s1", 3, "1,526,860.00", "0.00"),
    @("Neneng Hikmatul",          "Anisa_s1", "s1", 3, "419,354.00",   "0.00"),
    @("Dwi Gusti Anggraini",      "Anisa_s1", "s1", 2, "476,086.00",   "0.00"),
    @("Okky Al Bana",             "Anisa_s1", "s1", 1, "282,397.00",   "0.00"),
    @("Riska Rahmayanti",         "Anisa_s1", "s1", 1, "176,557.00",   "0.00"),
    @("Anisa Inraniwi",           "Anisa_s1", "s1", 1, "672,015.00",   "0.00"),
    @("Fachrul Rozi",             "Anisa_s1", "s1", 2, "359,327.00",   "0.00"),
    @("Ramesintia Sinaga",        "Anisa_s1", "s1", 2, "574,908.00",   "0.00"),
    @("Ayu Lestari",              "Anisa_s1", "s1", 2, "1,669,177.00", "0.00"),
    @("Dimas Kuat Anggowo",       "Anisa_s1", "s1", 5, "3,766,475.00", "0.00"),
    @("Gloriana Yesica",          "Anisa_s1", "s1", 2, "209,993.00",   "0.00")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $ws.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 0
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 0

    $r = $r + 1
}

# The old sheet had 14 data rows (rows 2-15); the new one only has 13
# (rows 2-14), so drop the now-unused last row entirely.
$ws.Rows(15).Delete()

# Column widths (best-fit sizing recorded by Excel after the edit).
$ws.Columns("A:A").ColumnWidth = 23.57
$ws.Columns("B:B").ColumnWidth = 8.71
$ws.Columns("C:C").ColumnWidth = 5.71
$ws.Columns("D:D").ColumnWidth = 22.14
$ws.Columns("E:E").ColumnWidth = 19.29
$ws.Columns("F:F").ColumnWidth = 16
$ws.Columns("G:G").ColumnWidth = 24.86
$ws.Columns("H:H").ColumnWidth = 9.71
$ws.Columns("I:I").ColumnWidth = 15.86
$ws.Columns("J:J").ColumnWidth = 27.14
$ws.Columns("K:K").ColumnWidth = 27.86
$ws.Columns("L:L").ColumnWidth = 26

$ws.Range("A1").Select()
